# Generate Report for Handoff
# - Bumps the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
#   for the six files that were just re-handed-off.
# - Sets the "Priority" column to "ht" for those same six rows on the
#   zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 10, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = Latest HO Xliff Generate Date
    $overview.Range("G$r").Value = "2016-08-29 18:23:31"

    # zh-cn sheet: column H = Latest Handoff Datetime, column E = Priority
    $zhcn.Range("H$r").Value = "2016-08-29 18:23:26"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: column H = Latest Handoff Datetime, column E = Priority
    $dede.Range("H$r").Value = "2016-08-29 18:23:31"
    $dede.Range("E$r").Value = "ht"
}
